# Add working set of sequences
# Fill columns G:N with "N/A" for the rows that currently only have A:F
# populated (these rows already hold "N/A" in C/D/F and originally had a
# shorter 1:6 span for columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,5,7,9,15,16,17,21,23,25,30,31,34,41,45,50,52,57,58,61,63,64,69,71,72,74,77,82,86,88,95,102,105,108,110,113,114,120,125,129,131,132,134,137,139,141,142,145,148,151,155,157,159,162,164,165,167,168,170,173,180,183,187,193)

foreach ($r in $rows) {
    $rng = $ws.Range("G" + $r + ":N" + $r)
    $rng.Value = "N/A"
}
